# Apply the commit "11 & 12 & 13" changes to the statistics workbook.
# The three affected sections are:
#   - 综合 (General) sheet: row 21 (全年民营经济增加值) value annotated with unit
#   - 贸易 (Trade) sheet: rows 1-3 & 11-13 (社会消费品零售总额/进出口总额 block) values updated
#   - 交通 (Transportation) sheet: rows 2, 5, 6 (铁路运营里程/铁路客运量/铁路货运量) values updated

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the cell to keep a literal text value (some strings such as
    # "+13%" would otherwise be auto-converted by Excel into a percentage
    # number). Formatting the cell as Text before assignment keeps the
    # value a plain string; clearing the formats afterwards removes the
    # now-unneeded "Text" number format so the cell keeps its original
    # (default) style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- 综合 (General) ---
$wsGeneral = $wb.Worksheets.Item("综合")
$wsGeneral.Range("C21").Value = "815.1亿元"

# --- 贸易 (Trade) ---
$wsTrade = $wb.Worksheets.Item("贸易")
$wsTrade.Range("C1").Value = "195亿元"
Set-TextValue $wsTrade.Range("D1") "+13%"
$wsTrade.Range("C2").Value = "29.9亿元"
Set-TextValue $wsTrade.Range("D2") "+26.2%"
$wsTrade.Range("C3").Value = "2013亿元"
Set-TextValue $wsTrade.Range("D3") "+8940.5%"
$wsTrade.Range("C11").Value = "3亿美元"
Set-TextValue $wsTrade.Range("D11") "-77.3%"
$wsTrade.Range("C12").Value = "33亿美元"
Set-TextValue $wsTrade.Range("D12") "-1%"
$wsTrade.Range("C13").Value = "51.1亿美元"
Set-TextValue $wsTrade.Range("D13") "+1734.7%"

# --- 交通 (Transportation) ---
$wsTransport = $wb.Worksheets.Item("交通")
$wsTransport.Range("C2").Value = "195公里"
$wsTransport.Range("C5").Value = "172.9万人"
$wsTransport.Range("C6").Value = "853万吨"
